# ReportingServer_Delavnica.pptx edit:
#  1) Slide 2 ("Agenda"): merge the two runs "Embedding R (Python) " + "code"
#     into a single run "Embedding R (Python) code".
#  2) Append a new slide (#8) "Give us your feedback" with a feedback
#     questionnaire in the content placeholder.

$p = $ppt.ActivePresentation

# --- 1) Fix the split run on slide 2 -----------------------------------
$agendaSlide = $p.Slides.Item(2)
$agendaBody  = $agendaSlide.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $agendaBody.Paragraphs().Count; $i++) {
    $candidate = $agendaBody.Paragraphs($i, 1)
    if ($candidate.Text -eq "Embedding R (Python) code") {
        # Re-assign the full paragraph text so the two runs collapse into one.
        $sub = $agendaBody.Characters($candidate.Start, $candidate.Length)
        $sub.Text = "Embedding R (Python) code"
        break
    }
}

# --- 2) Add the new "feedback" slide at the end -------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)  # ppLayoutText (Title and Content)

# Title placeholder
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Give us your feedback"
$title.LanguageID = "sl-SI"

# Content placeholder
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Give us your feedback"
$body.LanguageID = "sl-SI"

$lines = @(
    "How usefull / useless did you find this workshop?",
    "What would make it better?",
    "Was a particular topic left out?",
    "Would you recommend this workshop to your collegue?",
    "",
    "",
    "Email to: "
)
foreach ($line in $lines) {
    $r = $body.InsertAfter("`r" + $line)
    $r.LanguageID = "sl-SI"
}

# Shrink text on overflow (matches the authored slide's normAutofit bodyPr)
$newSlide.Shapes.Item(2).TextFrame.AutoSize = 2
